$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.5989378127028715
$ws.Range("C2").Value = 0.06512197661199792
$ws.Range("E2").Value = 0.09207146401279331
$ws.Range("F2").Value = 0.4443680307746263
$ws.Range("G2").Value = 1.145091058403949
$ws.Range("H2").Value = 1.13755260135764
$ws.Range("I2").Value = 1.123262273959412
$ws.Range("K2").Value = 0.3772425275440696
$ws.Range("L2").Value = 0.2144301975427965
$ws.Range("M2").Value = 0.1763691086489203
$ws.Range("N2").Value = 2.234211245682259
$ws.Range("B3").Value = 0.5665629001256605
$ws.Range("C3").Value = 0.06235449474846178
$ws.Range("E3").Value = 0.09225022174989128
$ws.Range("F3").Value = 0.387822817061874
$ws.Range("G3").Value = 1.150757544239781
$ws.Range("H3").Value = 1.144724194956851
$ws.Range("I3").Value = 1.131205857683689
$ws.Range("K3").Value = 0.3454086706624082
$ws.Range("L3").Value = 0.2118361181208499
$ws.Range("M3").Value = 0.1701809890785846
$ws.Range("N3").Value = 2.255917883416234
$ws.Range("B4").Value = 0.5469323596284141
$ws.Range("C4").Value = 0.0606318390217595
$ws.Range("E4").Value = 0.09239024742123547
$ws.Range("F4").Value = 0.3531389305168915
$ws.Range("G4").Value = 1.154859441009783
$ws.Range("H4").Value = 1.149570685602839
$ws.Range("I4").Value = 1.136570968257015
$ws.Range("K4").Value = 0.325984533796543
$ws.Range("L4").Value = 0.2103461761103915
$ws.Range("M4").Value = 0.1664636170782146
$ws.Range("N4").Value = 2.269926114499752
$ws.Range("B5").Value = 0.5389955528331427
$ws.Range("C5").Value = 0.05992394189422612
$ws.Range("E5").Value = 0.09245493871696908
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 1.156687521659464
$ws.Range("H5").Value = 1.151657167967308
$ws.Range("I5").Value = 1.138879972878446
$ws.Range("K5").Value = 0.3181000535234659
$ws.Range("L5").Value = 0.2097649353289057
$ws.Range("M5").Value = 0.1649695199489045
$ws.Range("N5").Value = 2.275805665856957
$ws.Range("B6").Value = 0.5376814593847143
$ws.Range("C6").Value = 0.05980603930315453
$ws.Range("E6").Value = 0.09246614202073289
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 1.157000525412606
$ws.Range("H6").Value = 1.152010363238716
$ws.Range("I6").Value = 1.139270791945105
$ws.Range("K6").Value = 0.31679272474328
$ws.Range("L6").Value = 0.209669988304455
$ws.Range("M6").Value = 0.1647226832687672
$ws.Range("N6").Value = 2.276792294668283
$ws.Range("B7").Value = 0.5468250661602951
$ws.Range("C7").Value = 0.06062231598431822
$ws.Range("E7").Value = 0.09239108895343051
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 1.15488346144592
$ws.Range("H7").Value = 1.149598373082014
$ws.Range("I7").Value = 1.136601611530885
$ws.Range("K7").Value = 0.3258780749156216
$ws.Range("L7").Value = 0.2103382322482332
$ws.Range("M7").Value = 0.1664433829603489
$ws.Range("N7").Value = 2.270004715443699
$ws.Range("B8").Value = 0.5877237972984517
$ws.Range("C8").Value = 0.06417259857299484
$ws.Range("E8").Value = 0.09212682868156463
$ws.Range("F8").Value = 0.4248636149813336
$ws.Range("G8").Value = 1.146915601863199
$ws.Range("H8").Value = 1.139933446031122
$ws.Range("I8").Value = 1.125900018063128
$ws.Range("K8").Value = 0.3662410906011928
$ws.Range("L8").Value = 0.2135144532415083
$ws.Range("M8").Value = 0.1742184523410515
$ws.Range("N8").Value = 2.241554521609633
$ws.Range("B9").Value = 0.6698751639467275
$ws.Range("C9").Value = 0.07094985097104711
$ws.Range("E9").Value = 0.09184800543126137
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 1.136234028966129
$ws.Range("H9").Value = 1.124493253796629
$ws.Range("I9").Value = 1.108782627802476
$ws.Range("K9").Value = 0.4463499205432697
$ws.Range("L9").Value = 0.2205566754393402
$ws.Range("M9").Value = 0.1901136624830571
$ws.Range("N9").Value = 2.191159441271571
$ws.Range("B10").Value = 0.7314038139879813
$ws.Range("C10").Value = 0.07581807244227434
$ws.Range("E10").Value = 0.0917881288845841
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 1.131405115537376
$ws.Range("H10").Value = 1.115286911346175
$ws.Range("I10").Value = 1.098563210377122
$ws.Range("K10").Value = 0.5057810188641838
$ws.Range("L10").Value = 0.2262243945203295
$ws.Range("M10").Value = 0.2021837093980636
$ws.Range("N10").Value = 2.15741959578139
$ws.Range("B11").Value = 0.759646266896965
$ws.Range("C11").Value = 0.07800902171375412
$ws.Range("E11").Value = 0.09179217053961786
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 1.129865051854509
$ws.Range("H11").Value = 1.111562088938967
$ws.Range("I11").Value = 1.094425596213561
$ws.Range("K11").Value = 0.5329412583388944
$ws.Range("L11").Value = 0.228909588850172
$ws.Range("M11").Value = 0.2077591347329459
$ws.Range("N11").Value = 2.142783030310971
$ws.Range("B12").Value = 0.7703768915062597
$ws.Range("C12").Value = 0.07883530009318918
$ws.Range("E12").Value = 0.09179818141188889
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 1.129376376615781
$ws.Range("H12").Value = 1.11021814034406
$ws.Range("I12").Value = 1.092932285933131
$ws.Range("K12").Value = 0.5432438112411262
$ws.Range("L12").Value = 0.229941724283691
$ws.Range("M12").Value = 0.2098824986483123
$ws.Range("N12").Value = 2.137342891463327
$ws.Range("B13").Value = 0.7680642747464788
$ws.Range("C13").Value = 0.0786574966998046
$ws.Range("E13").Value = 0.09179668787671602
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 1.129477416618684
$ws.Range("H13").Value = 1.110504624120665
$ws.Range("I13").Value = 1.093250627484466
$ws.Range("K13").Value = 0.5410241971732432
$ws.Range("L13").Value = 0.229718755533483
$ws.Range("M13").Value = 0.2094246590888389
$ws.Range("N13").Value = 2.138509967179742
$ws.Range("B14").Value = 0.7605283669570326
$ws.Range("C14").Value = 0.07807706794203284
$ws.Range("E14").Value = 0.09179257536473706
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 1.129822953539474
$ws.Range("H14").Value = 1.111450187757242
$ws.Range("I14").Value = 1.094301267353202
$ws.Range("K14").Value = 0.53378850580566
$ws.Range("L14").Value = 0.2289941967749058
$ws.Range("M14").Value = 0.2079335838539649
$ws.Range("N14").Value = 2.142333414071548
$ws.Range("B15").Value = 0.7559170551306522
$ws.Range("C15").Value = 0.07772109787065062
$ws.Range("E15").Value = 0.09179063928901599
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 1.130046916009164
$ws.Range("H15").Value = 1.112038039643082
$ws.Range("I15").Value = 1.094954388284847
$ws.Range("K15").Value = 0.5293587115079958
$ws.Range("L15").Value = 0.2285523755954415
$ws.Range("M15").Value = 0.20702182604375
$ws.Range("N15").Value = 2.144688724035188
$ws.Range("B16").Value = 0.7295631444985133
$ws.Range("C16").Value = 0.07567441424227184
$ws.Range("E16").Value = 0.09178849267599531
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 1.131518983181863
$ws.Range("H16").Value = 1.115539653912279
$ws.Range("I16").Value = 1.098843900603711
$ws.Range("K16").Value = 0.504008511099471
$ws.Range("L16").Value = 0.2260510567642768
$ws.Range("M16").Value = 0.2018210373578242
$ws.Range("N16").Value = 2.15839045807636
$ws.Range("B17").Value = 0.7134602378067427
$ws.Range("C17").Value = 0.07441279619401087
$ws.Range("E17").Value = 0.09179517475958399
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("G17").Value = 1.13259029248205
$ws.Range("H17").Value = 1.117806382218774
$ws.Range("I17").Value = 1.101360931692184
$ws.Range("K17").Value = 0.4884886636434373
$ws.Range("L17").Value = 0.2245439177666384
$ws.Range("M17").Value = 0.1986521385986961
$ws.Range("N17").Value = 2.166978384354838
$ws.Range("B18").Value = 0.7042220943304187
$ws.Range("C18").Value = 0.07368492211226396
$ws.Range("E18").Value = 0.09180196337038815
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 1.133268281261181
$ws.Range("H18").Value = 1.119153744564187
$ws.Range("I18").Value = 1.102856782795754
$ws.Range("K18").Value = 0.4795738210566753
$ws.Range("L18").Value = 0.2236871189970771
$ws.Range("M18").Value = 0.1968374507418105
$ws.Range("N18").Value = 2.171984933009066
$ws.Range("B19").Value = 0.7010983252608014
$ws.Range("C19").Value = 0.0734380940046151
$ws.Range("E19").Value = 0.09180476826463924
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("G19").Value = 1.133508448075261
$ws.Range("H19").Value = 1.119617428069702
$ws.Range("I19").Value = 1.103371517269359
$ws.Range("K19").Value = 0.4765574403569133
$ws.Range("L19").Value = 0.2233987527575181
$ws.Range("M19").Value = 0.1962244023017092
$ws.Range("N19").Value = 2.17369157042233
$ws.Range("B20").Value = 0.7151719579810276
$ws.Range("C20").Value = 0.07454732776005812
$ws.Range("E20").Value = 0.09179415874551111
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 1.132469853053934
$ws.Range("H20").Value = 1.117560572994535
$ws.Range("I20").Value = 1.101088009072519
$ws.Range("K20").Value = 0.4901395631841012
$ws.Range("L20").Value = 0.2247033137814469
$ws.Range("M20").Value = 0.1989886482814995
$ws.Range("N20").Value = 2.166057250301948
$ws.Range("B21").Value = 0.7627408777791516
$ws.Range("C21").Value = 0.07824764575966014
$ws.Range("E21").Value = 0.09179366185301063
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("G21").Value = 1.129718895180289
$ws.Range("H21").Value = 1.111170646710647
$ws.Range("I21").Value = 1.093990673814545
$ws.Range("K21").Value = 0.5359133301987242
$ws.Range("L21").Value = 0.2292066022341004
$ws.Range("M21").Value = 0.2083712220466509
$ws.Range("N21").Value = 2.141207595251862
$ws.Range("B22").Value = 0.794038462656971
$ws.Range("C22").Value = 0.08064628319495171
$ws.Range("E22").Value = 0.09181944407727727
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 1.128471924041818
$ws.Range("H22").Value = 1.107382407489425
$ws.Range("I22").Value = 1.08978065234313
$ws.Range("K22").Value = 0.565931308991992
$ws.Range("L22").Value = 0.2322389736134625
$ws.Range("M22").Value = 0.2145735870472194
$ws.Range("N22").Value = 2.12556386164902
$ws.Range("B23").Value = 0.7773154444557804
$ws.Range("C23").Value = 0.07936788723294796
$ws.Range("E23").Value = 0.0918033006366219
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 1.129087015137742
$ws.Range("H23").Value = 1.109368779589943
$ws.Range("I23").Value = 1.091988413232549
$ws.Range("K23").Value = 0.5499009243326611
$ws.Range("L23").Value = 0.2306123983321271
$ws.Range("M23").Value = 0.2112568705412698
$ws.Range("N23").Value = 2.13385857938165
$ws.Range("B24").Value = 0.7143980280986284
$ws.Range("C24").Value = 0.07448651399442952
$ws.Range("E24").Value = 0.0917946089048538
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 1.132524110353231
$ws.Range("H24").Value = 1.117671565689434
$ws.Range("I24").Value = 1.101211245451481
$ws.Range("K24").Value = 0.4893931674130272
$ws.Range("L24").Value = 0.2246312207009993
$ws.Range("M24").Value = 0.1988364899612378
$ws.Range("N24").Value = 2.166473479341127
$ws.Range("B25").Value = 0.6474440528338619
$ws.Range("C25").Value = 0.06913604150216202
$ws.Range("E25").Value = 0.09189790594840197
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 1.138593916943606
$ws.Range("H25").Value = 1.128294579852778
$ws.Range("I25").Value = 1.112999359795438
$ws.Range("K25").Value = 0.4245767571747479
$ws.Range("L25").Value = 0.2185646902249871
$ws.Range("M25").Value = 0.1857445309808128
$ws.Range("N25").Value = 2.204215420786491
